$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sat on the empty page-break
#    paragraph right before "Changement hardware du serveur".
# ---------------------------------------------------------------------------
$xmlPageBreak = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="63FDD114" w14:textId="77777777" w:rsidR="005928F6" w:rsidRDefault="005928F6"><w:pPr><w:jc w:val="left"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:sz w:val="36"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:sz w:val="36"/><w:szCs w:val="32"/></w:rPr><w:br w:type="page"/></w:r></w:p>
'@
$d.Paragraphs(29).Range.InsertXML($xmlPageBreak)

# ---------------------------------------------------------------------------
# 2) "Changement hardware du serveur" heading - bookmark id 2 -> 1
# ---------------------------------------------------------------------------
$xmlChangement = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="67E0F23C" w14:textId="16ED0032" w:rsidR="0010668A" w:rsidRDefault="00EA2072" w:rsidP="002D6EF9"><w:pPr><w:pStyle w:val="Titre2"/></w:pPr><w:bookmarkStart w:id="1" w:name="_Toc126657547"/><w:r><w:lastRenderedPageBreak/><w:t>Changement</w:t></w:r><w:r w:rsidR="002D6EF9"><w:t xml:space="preserve"> hardware du serveur</w:t></w:r><w:bookmarkEnd w:id="1"/></w:p>
'@
$d.Paragraphs(30).Range.InsertXML($xmlChangement)

# ---------------------------------------------------------------------------
# 3) "Choix du nouveau matériel" / "Justification financière" headings -
#    bookmark ids 3 -> 2 and 4 -> 3.
# ---------------------------------------------------------------------------
$xmlChoixJustif = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="65E586BC" w14:textId="64221762" w:rsidR="00EA2072" w:rsidRDefault="00EA2072" w:rsidP="00EA2072"><w:pPr><w:pStyle w:val="Titre3"/></w:pPr><w:bookmarkStart w:id="2" w:name="_Toc126657548"/><w:r><w:t>Choix du nouveau matériel</w:t></w:r><w:bookmarkEnd w:id="2"/></w:p><w:p w14:paraId="0AE7EFC2" w14:textId="687852C9" w:rsidR="00EA2072" w:rsidRPr="00EA2072" w:rsidRDefault="00EA2072" w:rsidP="00EA2072"><w:pPr><w:pStyle w:val="Titre3"/></w:pPr><w:bookmarkStart w:id="3" w:name="_Toc126657549"/><w:r><w:t>Justification financière</w:t></w:r><w:bookmarkEnd w:id="3"/></w:p>
'@
$rngChoixJustif = $d.Range($d.Paragraphs(34).Range.Start, $d.Paragraphs(35).Range.End)
$rngChoixJustif.InsertXML($xmlChoixJustif)

# ---------------------------------------------------------------------------
# 4) "Migration du système d'exploitation du serveur" heading -
#    bookmark id 5 -> 4.
# ---------------------------------------------------------------------------
$xmlMigrationOs = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="50701695" w14:textId="3D623DF4" w:rsidR="003B3046" w:rsidRDefault="003B3046" w:rsidP="003B3046"><w:pPr><w:pStyle w:val="Titre2"/></w:pPr><w:bookmarkStart w:id="4" w:name="_Toc126657550"/><w:r><w:t>Migration du système d’exploitation du serveur</w:t></w:r><w:bookmarkEnd w:id="4"/></w:p>
'@
$d.Paragraphs(36).Range.InsertXML($xmlMigrationOs)

# ---------------------------------------------------------------------------
# 5) "OS actuel : Microsoft Windows Server 2003..." paragraph - add a
#    French typography space before the colon (grammar-checker fix),
#    move the _GoBack bookmark here, and wrap the fix in gramStart/gramEnd
#    proof markers.
# ---------------------------------------------------------------------------
$xmlOsActuel = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="50D46902" w14:textId="54C4B6F6" w:rsidR="00B80752" w:rsidRDefault="00B80752" w:rsidP="00B80752"><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr><w:r w:rsidRPr="00D934E5"><w:rPr><w:b/><w:lang w:val="de-CH"/></w:rPr><w:t>OS</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:lang w:val="de-CH"/></w:rPr><w:t>actuel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00207D0E"><w:rPr><w:b/><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="5" w:name="_GoBack"/><w:bookmarkEnd w:id="5"/><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="002572C9"><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> Microsoft Windows Server 2003 Standard E</w:t></w:r><w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>dition</w:t></w:r></w:p>
'@
$d.Paragraphs(37).Range.InsertXML($xmlOsActuel)

Write-Output "edits applied"
